$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-18 (Name, Position, Team)
$data = @(
  @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
  @("Jamal Murray", "PG,SG", "Denver Nuggets"),
  @("Luguentz Dort", "SG,SF", "Oklahoma City Thunder"),
  @("Victor Wembanyama", "C", "San Antonio Spurs"),
  @("Malik Beasley", "SG", "Detroit Pistons"),
  @("Myles Turner", "C", "Indiana Pacers"),
  @("Josh Hart", "SF,PF", "New York Knicks"),
  @("Domantas Sabonis", "C", "Sacramento Kings"),
  @("Tari Eason", "SF,PF", "Houston Rockets"),
  @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
  @("Ochai Agbaji", "SG,SF", "Toronto Raptors"),
  @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
  @("Alexandre Sarr", "PF,C", "Washington Wizards"),
  @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
  @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
  @("Ja Morant", "PG", "Memphis Grizzlies"),
  @("Bradley Beal", "PG,SG,SF", "Phoenix Suns")
)

$row = 2
foreach ($entry in $data) {
  $ws.Cells.Item($row, 1).Value = $entry[0]
  $ws.Cells.Item($row, 2).Value = $entry[1]
  $ws.Cells.Item($row, 3).Value = $entry[2]
  $row = $row + 1
}

# The old sheet had one extra row (row 19) that no longer exists; remove it.
$ws.Range("A19:C19").Delete()
